# Add analyzer (ID/tag) values to the asset rows in the "Analyzer" section
# of the balance-sheet / income-statement table on Sheet1, column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row -> analyzer tag text (column A)
$tags = @{
    33 = "ts无此数据"
    34 = "T004 - 仅固定资产"
    36 = "ts无此科目"
    37 = "ts无此数据"
    39 = "T004"
    41 = "T004"
    43 = "ts无此数据"
    44 = "ts无此数据"
    50 = "难以自动判定"
    51 = "T004 - 仅固定资产"
    52 = "难以自动判定"
    60 = "T005"
    61 = "T005"
    65 = "T005"
    67 = "T005"
    69 = "T005"
    71 = "T005"
    72 = "T005"
    73 = "T005"
    74 = "T005"
}

foreach ($row in $tags.Keys) {
    $ws.Cells.Item($row, 1).Value = $tags[$row]
}

# Scroll the frozen pane down and move the active selection, matching the
# author's final view position in the sheet.
$ws.Activate()
$ws.Cells.Item(74, 1).Select()
$excel.ActiveWindow.ScrollRow = 65

